$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.507.25"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.220.36"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'270.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("D6").Value = "'92.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.07%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").Value = "'45.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.70%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "'8.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.42%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "2.555.70"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "2.228.26"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'0.800"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "43.488.87"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'6.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "'70.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "'232.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "'9.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'11.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.66%  "
$ws.Range("D27").Value = "'2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.18%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'41.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").Value = "'3.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "'172.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "'0.0920"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.78%  "
$ws.Range("D33").Value = "'20.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  -4.98%  "
$ws.Range("D39").Value = "'3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.74%  "
$ws.Range("D40").Value = "'12.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("E41").Value = "  +9.30%  "
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").Value = "'63.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").Value = "'0.0987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").Value = "'100.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "

Write-Output "Applied cryptos update"
